$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 388; $r++) {
    $ws.Range("C$r").Value = 46062
}
